# Update Gip-Gipr NATMI output worksheet with newly-computed TPM values.
# Rows 4-7 (the MuSCs-related sending/target cluster combinations) are no
# longer present in the new TPM run, so they are deleted. The remaining
# two rows (ECs->ECs and ECs->FAPs) keep their identifiers but get
# refreshed numeric statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete MuSCs-related rows (old rows 4 through 7).
$ws.Rows("4:7").Delete()

# Refresh the numeric columns (G:T) for the two remaining rows with the
# newly computed TPM-derived values.

# Row 2: ECs (sending) -> Gip -> Gipr -> ECs (target)
$ws.Range("G2").Value = 0.01348466666666667
$ws.Range("H2").Value = 0.040454
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.03822966666666667
$ws.Range("N2").Value = 0.114689
$ws.Range("O2").Value = 0.287343161228253
$ws.Range("P2").Value = 0.287343161228253
$ws.Range("Q2").Value = 0.0005155143117777778
$ws.Range("R2").Value = 0.004639628805999999
$ws.Range("S2").Value = 0.287343161228253
$ws.Range("T2").Value = 0.287343161228253

# Row 3: ECs (sending) -> Gip -> Gipr -> FAPs (target)
$ws.Range("G3").Value = 0.01348466666666667
$ws.Range("H3").Value = 0.040454
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.09481566666666667
$ws.Range("N3").Value = 0.284447
$ws.Range("O3").Value = 0.7126568387717469
$ws.Range("P3").Value = 0.712656838771747
$ws.Range("Q3").Value = 0.001278557659777778
$ws.Range("R3").Value = 0.011507018938
$ws.Range("S3").Value = 0.7126568387717469
$ws.Range("T3").Value = 0.712656838771747
